# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# for the 48511851-2b57-4eea-9dfd-fa2a429c2fad.md file across all three report sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: row 7 ("48511851-2b57-4eea-9dfd-fa2a429c2fad.md"),
# column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-09-02 06:48:54"

# zh-cn sheet: row 7, column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-09-02 06:48:49"

# de-de sheet: row 7, column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-09-02 06:48:54"
